$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "28.180.79"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.834.56"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.07"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5121"
$ws.Range("E7").Value = "  -2.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3962"
$ws.Range("E8").Value = "  +3.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09973"
$ws.Range("E9").Value = "  +24.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.113"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.08"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.492"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.71"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.424"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "1.820.11"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001142"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.00"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06629"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.42"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9988"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.068"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "28.210.38"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.252"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.461"
$ws.Range("E26").Value = "  +4.49%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.85"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.37"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").Value = "2.044.11"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.29"
$ws.Range("E30").Value = "  +4.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1098"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.065"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.642"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.639"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06916"
$ws.Range("E35").Value = "  -4.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.176"
$ws.Range("E36").Value = "  +6.58%  "
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2181"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.64"
$ws.Range("E39").Value = "  -5.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.028"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6295"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9975"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.156"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.34"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6014"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.294"
$ws.Range("E46").Value = "  -5.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.710"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.29"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.994"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.193"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06799"
$ws.Range("E51").Value = "  -0.27%  "
